# Non-Technical.pptx edit:
#  1. Slide 1 ("King County Homes" title slide) - set the (previously
#     empty) title placeholder text.
#  2. Append a brand-new slide (sldId 257) using the "Title and Content"
#     layout (the 2nd layout of the slide master), left with its default
#     empty placeholders - matching the commit "didn't really do anything".

$p = $ppt.ActivePresentation

# --- 1. Title text on the existing (first) slide ------------------------
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "King County Homes"

# --- 2. New slide appended at the end, "Title and Content" layout -------
$s2 = $p.Slides.Add($p.Slides.Count + 1, 2)
